$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.079.65'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.214.85'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.96'
$ws.Range('E5').Value = '  -2.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.625'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.21'
$ws.Range('E7').Value = '  -1.41%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.64'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0953'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.09'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.548.27'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.17'
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.838'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.224.48'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.940.55'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  +9.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.58'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.22'
$ws.Range('E22').Value = '  +16.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.40'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  -7.12%  '
$ws.Range('E25').Value = '  +3.51%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.13'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.57'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  +9.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0790'
$ws.Range('E33').Value = '  -2.98%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '28.76'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('E36').Value = '  -7.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.26'
$ws.Range('E37').Value = '  -4.66%  '
$ws.Range('E38').Value = '  -3.01%  '
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('E40').Value = '  -2.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '64.85'
$ws.Range('E41').Value = '  +4.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.61'
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('E43').Value = '  -2.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.72'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.95'
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  +5.98%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.422.69'
$ws.Range('E51').Value = '  -2.07%  '
